# Apply latest cryptocurrency price/volume data to the worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a "Price"-column cell while forcing it to remain
# a text string (the sheet stores prices like "27.919.46" / "7.60" as text,
# not numbers). Without this, Excel would silently reinterpret the text as a
# number (and e.g. drop the trailing zero in "7.60").
function Set-PriceText($addr, $val) {
    $c = $ws.Range($addr)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

# --- Price column (D) updates ---
Set-PriceText "D2" "27.919.46"
Set-PriceText "D3" "1.638.97"
Set-PriceText "D5" "213.59"
Set-PriceText "D8" "23.61"
Set-PriceText "D9" "0.262"
Set-PriceText "D12" "1.871.48"
Set-PriceText "D13" "1.648.54"
Set-PriceText "D17" "27.909.36"
Set-PriceText "D18" "231.73"
Set-PriceText "D20" "7.60"
Set-PriceText "D22" "10.84"
Set-PriceText "D25" "151.77"
Set-PriceText "D27" "15.75"
Set-PriceText "D34" "1.412.80"
Set-PriceText "D40" "0.914"
Set-PriceText "D41" "1.02"
Set-PriceText "D43" "66.27"
Set-PriceText "D44" "1.83"
Set-PriceText "D46" "2.20"
Set-PriceText "D47" "1.780.27"
Set-PriceText "D48" "88.31"
Set-PriceText "D50" "0.0506"
Set-PriceText "D51" "7.66"

# --- Coin name / link / volume column (B, C, E) updates ---
$ws.Range("E2").Value = "  +1.41%  "
$ws.Range("E3").Value = "  +0.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("E5").Value = "  +0.74%  "
$ws.Range("E6").Value = "  +0.65%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +1.09%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  +0.62%  "
$ws.Range("E11").Value = "  -0.33%  "
$ws.Range("E12").Value = "  +0.77%  "
$ws.Range("E13").Value = "  +1.27%  "
$ws.Range("E14").Value = "  +1.15%  "
$ws.Range("E15").Value = "  +4.02%  "
$ws.Range("E16").Value = "  +0.99%  "
$ws.Range("E17").Value = "  +1.44%  "
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("E19").Value = "  +0.95%  "
$ws.Range("E20").Value = "  +0.35%  "
$ws.Range("E21").Value = "  +0.02%  "
$ws.Range("E22").Value = "  +3.43%  "
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("E24").Value = "  -3.75%  "
$ws.Range("E25").Value = "  +1.78%  "
$ws.Range("E26").Value = "  +0.52%  "
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("E28").Value = "  +0.09%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  +0.94%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("E32").Value = "  +1.95%  "
$ws.Range("E33").Value = "  +1.36%  "
$ws.Range("E34").Value = "  -3.77%  "
$ws.Range("E35").Value = "  +1.16%  "
$ws.Range("E36").Value = "  +0.61%  "
$ws.Range("E37").Value = "  +1.79%  "
$ws.Range("E38").Value = "  +0.12%  "
$ws.Range("E39").Value = "  +0.16%  "
$ws.Range("E40").Value = "  -2.84%  "
$ws.Range("E41").Value = "  -1.31%  "
$ws.Range("E42").Value = "  +0.05%  "
$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E43").Value = "  -2.20%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E44").Value = "  +3.75%  "
$ws.Range("E45").Value = "  +1.52%  "
$ws.Range("E46").Value = "  +0.15%  "
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("E48").Value = "  +1.09%  "
$ws.Range("E49").Value = "  +0.94%  "
$ws.Range("E50").Value = "  +0.42%  "
$ws.Range("E51").Value = "  -0.45%  "
